# Update the "Ready for handoff" status text to "In Translation" across all
# worksheets that contain it, then re-run AutoFit on the affected columns so
# the column widths shrink to match the new (shorter) text, matching the
# behaviour Excel performs automatically when AutoFit is applied.

$wb = $excel.ActiveWorkbook

$oldText = "Ready for handoff"
$newText = "In Translation"

# --- Overview sheet: columns E and F hold the status for each locale -------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newText
$wsOverview.Range("F2").Value = $newText
$wsOverview.Range("E3").Value = $newText
$wsOverview.Range("F3").Value = $newText
$wsOverview.Range("E2:F3").Columns.AutoFit() | Out-Null

# --- zh-cn sheet: column C holds Status -------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newText
$wsZh.Range("C3").Value = $newText
$wsZh.Range("C2:C3").Columns.AutoFit() | Out-Null

# --- de-de sheet: column C holds Status -------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newText
$wsDe.Range("C3").Value = $newText
$wsDe.Range("C2:C3").Columns.AutoFit() | Out-Null
